# Updated cryptos list on Fri Apr 14 10:40:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2  = @{ D = "30.816.93";   E = "  +1.94%  " }
    3  = @{ D = "2.114.12";    E = "  +6.40%  " }
    4  = @{                    E = "  -0.01%  " }
    5  = @{ D = "333.02";      E = "  +2.99%  " }
    6  = @{ D = "1.000";       E = "  +0.02%  " }
    7  = @{ D = "0.5326";      E = "  +4.43%  " }
    8  = @{ D = "0.4380";      E = "  +6.64%  " }
    9  = @{ D = "0.08987";     E = "  +1.84%  " }
    10 = @{ D = "47.12";       E = "  +11.06%  " }
    11 = @{ D = "1.181";       E = "  +4.66%  " }
    12 = @{ D = "25.01";       E = "  +3.48%  " }
    13 = @{ D = "2.118.97";    E = "  +6.98%  " }
    14 = @{ D = "6.764";       E = "  +4.35%  " }
    15 = @{ D = "7.825";       E = "  +5.88%  " }
    16 = @{ D = "97.07";       E = "  +3.38%  " }
    17 = @{                    E = "  -0.36%  " }
    18 = @{ D = "0.00001132";  E = "  +1.09%  " }
    19 = @{ D = "0.06669";     E = "  +2.01%  " }
    20 = @{ D = "19.14";       E = "  +1.88%  " }
    21 = @{ D = "0.9997";      E = "  -0.09%  " }
    22 = @{ D = "6.342";       E = "  +4.45%  " }
    23 = @{ D = "30.886.84";   E = "  +1.98%  " }
    24 = @{ D = "12.37";       E = "  +7.70%  " }
    25 = @{ D = "2.366.03";    E = "  +6.95%  " }
    26 = @{ D = "2.270";       E = "  +2.65%  " }
    27 = @{ D = "22.79";       E = "  +1.34%  " }
    28 = @{ D = "2.581";       E = "  +9.02%  " }
    29 = @{ D = "163.31";      E = "  +0.35%  " }
    30 = @{ D = "133.49";      E = "  +2.07%  " }
    31 = @{ D = "1.181";       E = "  +4.09%  " }
    32 = @{ D = "0.1082";      E = "  +2.72%  " }
    33 = @{ D = "6.251";       E = "  +3.56%  " }
    34 = @{ D = "4.016";       E = "  +5.73%  " }
    35 = @{ D = "1.563";       E = "  +19.12%  " }
    36 = @{ D = "0.02607";     E = "  +5.13%  " }
    37 = @{ D = "12.91";       E = "  +10.09%  " }
    38 = @{ D = "5.545";       E = "  +3.03%  " }
    39 = @{ D = "0.06761";     E = "  +4.20%  " }
    40 = @{ D = "9.484";       E = "  +6.51%  " }
    41 = @{                    E = "  +5.26%  " }
    42 = @{ D = "0.6870";      E = "  +4.79%  " }
    43 = @{ D = "1.252";       E = "  +2.63%  " }
    44 = @{ D = "0.6470";      E = "  +5.67%  " }
    45 = @{ D = "14.18";       E = "  +4.63%  " }
    46 = @{ D = "0.9995";      E = "  -0.02%  " }
    47 = @{ D = "2.228";       E = "  +2.01%  " }
    48 = @{ D = "3.670";       E = "  +0.41%  " }
    49 = @{ D = "1.276";       E = "  +4.61%  " }
    50 = @{ D = "83.09";       E = "  +4.72%  " }
    51 = @{ D = "121.86";      E = "  -1.49%  " }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    if ($vals.ContainsKey("D")) {
        $cellD = $ws.Cells.Item($r, 4)
        # Force the cell to remain plain text (matches original inline-string
        # cells), so values such as "333.02" or "1.000" are not re-interpreted
        # as numbers by Excel's automatic type detection.
        $cellD.NumberFormat = "@"
        $cellD.Value = $vals["D"]
        $cellD.Style = "Normal"
    }
    $ws.Cells.Item($r, 5).Value = $vals["E"]
}
